# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.681.05'
$ws.Range('E2').Value = '  +0.66%  '
$ws.Range('D3').Value = '2.278.74'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '504.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.70'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.35%  '
$ws.Range('E7').Value = '  -0.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.528'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '2.294.88'
$ws.Range('E9').Value = '  +0.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0967'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.99%  '
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('E12').Value = '  +3.44%  '
$ws.Range('E13').Value = '  +3.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.40%  '
$ws.Range('D15').Value = '2.685.76'
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('D16').Value = '54.744.40'
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').Value = '2.308.81'
$ws.Range('E18').Value = '  +0.52%  '
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.13'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '306.97'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.10%  '
$ws.Range('E25').Value = '  -0.49%  '
$ws.Range('E26').Value = '  +0.03%  '
$ws.Range('E27').Value = '  +1.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '171.48'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.63'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0701'
$ws.Range('E30').Value = '  +2.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.04'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.85%  '
$ws.Range('E32').Value = '  +2.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.94'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.26%  '
$ws.Range('E35').Value = '  -0.14%  '
$ws.Range('E36').Value = '  -0.77%  '
$ws.Range('E37').Value = '  -0.20%  '
$ws.Range('E38').Value = '  +1.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.40'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.68%  '
$ws.Range('E40').Value = '  +0.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.42'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.05'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.32%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '126.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '251.79'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.00%  '
$ws.Range('E46').Value = '  +1.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0901'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('E48').Value = '  +0.36%  '
$ws.Range('E49').Value = '  +0.79%  '
$ws.Range('E50').Value = '  +0.47%  '
